$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "In Translation" for every cell that carries that
# status text across the Overview / zh-cn / de-de sheets. Changing the text
# also shrinks the Status columns so they get re-autofitted to the new,
# narrower content below.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# Narrow the affected status columns to match the new, shorter text
# (~13.41 characters wide).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
